$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.006876353814593728
$ws.Range("C2").Value = 3486408257.158922
$ws.Range("D2").Value = 1935279062.313128
$ws.Range("E2").Value = 5548678842208.939
$ws.Range("G2").Value = 5554100529528.418
